$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set cell values for new test routine rows (mirsum / mirpeaks), rows 151-179
$ws.Range("B151").Value = "mirsum"
$ws.Range("D151").Value = "sig.sum"
$ws.Range("C152").Value = "…"
$ws.Range("E152").Value = "…"
$ws.Range("C153").Value = "…, ‘Center’"
$ws.Range("E153").Value = "…, ‘Center’"
$ws.Range("C154").Value = "…, ‘Mean’"
$ws.Range("E154").Value = "…, ‘Mean’"
$ws.Range("B155").Value = "mirpeaks"
$ws.Range("D155").Value = "sig.peaks"
$ws.Range("C156").Value = "…"
$ws.Range("E156").Value = "…"
$ws.Range("C157").Value = "..., 'Total', m"
$ws.Range("E157").Value = "..., 'Total', m"
$ws.Range("C158").Value = "..., 'Total', m, 'NoBegin'"
$ws.Range("E158").Value = "..., 'Total', m, 'NoBegin'"
$ws.Range("C159").Value = "..., 'Total', m, 'NoEnd'"
$ws.Range("E159").Value = "..., 'Total', m, 'NoEnd'"
$ws.Range("C160").Value = "..., 'Order', 'Amplitude'"
$ws.Range("E160").Value = "..., 'Order', 'Amplitude'"
$ws.Range("C161").Value = "..., 'Order', 'Abscissa'"
$ws.Range("E161").Value = "..., 'Order', 'Abscissa'"
$ws.Range("C162").Value = "..., 'Valleys'"
$ws.Range("E162").Value = "..., 'Valleys'"
$ws.Range("C163").Value = "..., 'Contrast', cthr"
$ws.Range("E163").Value = "..., 'Contrast', cthr"
$ws.Range("C164").Value = "..., 'SelectFirst', fthr"
$ws.Range("E164").Value = "..., 'SelectFirst', fthr"
$ws.Range("C165").Value = "..., 'Contrast', cthr, 'SelectFirst', fthr"
$ws.Range("E165").Value = "..., 'Contrast', cthr, 'SelectFirst', fthr"
$ws.Range("C166").Value = "..., 'Threshold', thr"
$ws.Range("E166").Value = "..., 'Threshold', thr"
$ws.Range("C167").Value = "..., 'Valleys'. 'Threshold', thr"
$ws.Range("E167").Value = "..., 'Valleys'. 'Threshold', thr"
$ws.Range("C168").Value = "..., 'Interpol', 'no'"
$ws.Range("E168").Value = "..., 'Interpol', 'no'"
$ws.Range("C169").Value = "..., 'Interpol', 'Quadratic'"
$ws.Range("E169").Value = "..., 'Interpol', 'Quadratic'"
$ws.Range("C170").Value = "..., 'Reso', r "
$ws.Range("E170").Value = "..., 'Reso', r "
$ws.Range("C171").Value = "..., 'Reso', r, 'First'"
$ws.Range("E171").Value = "..., 'Reso', r, 'First'"
$ws.Range("C172").Value = "…, ‘Pref’, c, std"
$ws.Range("F172").Value = "NI"
$ws.Range("C173").Value = "…, ‘Nearest’, t, s"
$ws.Range("F173").Value = "NI"
$ws.Range("C174").Value = "..., 'Normalize', 'Global' "
$ws.Range("E174").Value = "..., 'Normalize', 'Global' "
$ws.Range("C175").Value = "..., 'Normalize', 'Local' "
$ws.Range("E175").Value = "..., 'Normalize', 'Local' "
$ws.Range("C176").Value = "…, ‘Extract’"
$ws.Range("F176").Value = "NI"
$ws.Range("C177").Value = "…,’Only’"
$ws.Range("F177").Value = "NI"
$ws.Range("C178").Value = "…, ‘Track’, t"
$ws.Range("F178").Value = "NI"
$ws.Range("C179").Value = "…, ‘CollapseTracks’, t"
$ws.Range("F179").Value = "NI"

# Apply empty-but-styled cells (E172/E173 carry style only, no text)
$ws.Range("E172").Value = ""
$ws.Range("E173").Value = ""

# Row heights matching the authored rows (sheet default is already 15,
# so only the rows whose height differs from that default need to be set)
$ws.Rows.Item(152).RowHeight = 13.8
$ws.Rows.Item(153).RowHeight = 13.8
$ws.Rows.Item(154).RowHeight = 13.8
$ws.Rows.Item(157).RowHeight = 13.8
$ws.Rows.Item(158).RowHeight = 13.8
$ws.Rows.Item(159).RowHeight = 13.8
$ws.Rows.Item(160).RowHeight = 13.8
$ws.Rows.Item(161).RowHeight = 13.8
$ws.Rows.Item(162).RowHeight = 13.8
$ws.Rows.Item(163).RowHeight = 13.8
$ws.Rows.Item(164).RowHeight = 13.8
$ws.Rows.Item(165).RowHeight = 13.8
$ws.Rows.Item(166).RowHeight = 13.8
$ws.Rows.Item(167).RowHeight = 13.8
$ws.Rows.Item(168).RowHeight = 13.8
$ws.Rows.Item(169).RowHeight = 13.8
$ws.Rows.Item(170).RowHeight = 13.8
$ws.Rows.Item(171).RowHeight = 13.8
$ws.Rows.Item(172).RowHeight = 13.8
$ws.Rows.Item(173).RowHeight = 13.8
$ws.Rows.Item(174).RowHeight = 13.8
$ws.Rows.Item(175).RowHeight = 13.8

# Remove the stray trailing formatting-only row left over at the bottom of the sheet
$ws.Rows.Item(1048576).Delete()

# Restore/update the active selection to match where the new content was added
$excel.ActiveWindow.ScrollRow = 153
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F171").Select()
